$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -2
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = -2
$ws.Range("F5").Value = -8
$ws.Range("F6").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("F10").Value = 10
$ws.Range("F11").Value = 4
$ws.Range("F12").Value = 6
